# Updates the hashcode values (column B) for the rows identified by their
# key in column A, matching the "Actualización automática hashcode" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Key = "05-050305TC"; NewValue = "a0c1161837786ed577bd398f7504ad26" },
    @{ Key = "05-050305TP"; NewValue = "dbb17ca4b52a4c7e5e94472e9b66584d" },
    @{ Key = "05-050104A";  NewValue = "68439b0181d7876541c13259662e62d3" },
    @{ Key = "05-050101A";  NewValue = "0dc061740719f94d60c3f2fb1a76b472" },
    @{ Key = "05-050102A";  NewValue = "aa9b18f3904e71ef4dadf88111858b4d" },
    @{ Key = "05-050301TP"; NewValue = "89e31980121a03ecb4d90a72f238e8a1" },
    @{ Key = "05-050309A";  NewValue = "cb2b48530b102a7818d954df99d33a88" },
    @{ Key = "05-050007TC"; NewValue = "7883f0f152cc9d9bb5a1fc710f211227" },
    @{ Key = "05-050007TP"; NewValue = "86c3466b53645a70143a60d23010a457" },
    @{ Key = "05-050305A";  NewValue = "bfd43c2f789ae217aee9d6a0c58b3db0" },
    @{ Key = "05-050306A";  NewValue = "d390d1e05d7bb974a4a42141255c0a4f" },
    @{ Key = "05-050304A";  NewValue = "32971f9e01b7e44aa184d2c517c589a3" },
    @{ Key = "05-050101TP"; NewValue = "f628fff06e904e2e47120b72d229abd7" },
    @{ Key = "05-050310TC"; NewValue = "d9e41eccb1727d9b81e0c2d1587a1e06" },
    @{ Key = "05-050005TP"; NewValue = "11352530e667e1d92b0f0b73ab121312" },
    @{ Key = "05-050006A";  NewValue = "b526e2e952a95b9a09ec2a8738f95769" },
    @{ Key = "05-050007A";  NewValue = "5cbb749084cfb11e073fabbd9fa5cca4" },
    @{ Key = "05-050308TP"; NewValue = "c3f60ea1fa19ab1c30e5690afe2c4a50" },
    @{ Key = "05-050005A";  NewValue = "320e9ebd681ed0347b22b3f3e81e84de" },
    @{ Key = "05-050304TP"; NewValue = "7068eaeabb596cde9800331635f8126a" },
    @{ Key = "05-050102TP"; NewValue = "a674c1abc8131bd1104e7863c9f31dd5" },
    @{ Key = "05-050006TC"; NewValue = "19e459ae140fd3ca9c68c0372a062362" },
    @{ Key = "05-050006TP"; NewValue = "ce02acf55c77ea096712c1a555e3035c" },
    @{ Key = "05-050104TC"; NewValue = "6b89c2b53a18291ef70de72b5ff36fa1" },
    @{ Key = "05-050104TM"; NewValue = "3b90ab400a44cba436858271a190263b" },
    @{ Key = "05-050104TP"; NewValue = "462b9661f05db7b33cc099b42a4fe747" },
    @{ Key = "05-050309TP"; NewValue = "73ac72d57a94466bf0648eef63be2fea" },
    @{ Key = "03-030032A";  NewValue = "c9c849f03081bb7a17b5eba5feebb7ea" }
)

$colA = $ws.Columns.Item(1)

foreach ($u in $updates) {
    $found = $colA.Find($u.Key, [Type]::Missing, [Type]::Missing, 1)
    if ($found -eq $null) {
        throw "Key not found: $($u.Key)"
    }
    $row = $found.Row
    $ws.Cells.Item($row, 2).Value = $u.NewValue
}
